# seller in add listing changes purpose
#
# Fill in the newly added "add listing" seller fields (AF2:AM2) on the
# first worksheet, then leave the view scrolled to the right with the
# last-entered cell selected - matching where the user ended up after
# typing the values in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AF2").Value = "gh"
$ws.Range("AG2").Value = "hfgh"
$ws.Range("AH2").Value = "fghfgh"
$ws.Range("AI2").Value = "fhdfg"
$ws.Range("AJ2").Value = "fgd"
$ws.Range("AK2").Value = "fgh"
$ws.Range("AL2").Value = "hfdgh"
$ws.Range("AM2").Value = "fg"

# Scroll the window over toward the new columns and select the final cell.
$excel.ActiveWindow.ScrollColumn = 30
$ws.Range("AM2").Select()
